# Apply the "change language to ind" edit to the santri (student) sheet.
# Strategy:
#  1. Remove the old per-column custom widths (<cols>) by deleting the
#     originally-customized columns (A:K) as a single EntireColumn.Delete()
#     call. This also wipes the old cell data/shared strings in one go.
#  2. Re-enter the new Indonesian-language headers and the two data rows,
#     cell by cell, in left-to-right / top-to-bottom order so that the
#     shared-strings table is rebuilt with the exact same de-duplicated
#     ordering as the target workbook.
#  3. A few cells (No Induk/NISN/Anak ke-) are numeric, so those are
#     written as real numbers instead of text.
#  4. Some text values look like ISO dates ("2020-08-04" etc.); those are
#     forced to remain plain text (instead of being auto-converted to a
#     date serial number) by temporarily marking the cell as Text format,
#     then clearing that formatting again once the literal string value
#     has been stored.
#  5. Finally the sheet selection is updated to match the saved state
#     (whole column A selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Strip the old custom column widths -------------------------------
$ws.Range("A1:K1").EntireColumn.Delete()

# --- helper for cells that must stay literal text even though they look
#     like dates -----------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- 2. Row 1: headers -----------------------------------------------------
$ws.Range("A1").Value = "Nama Lengkap"
$ws.Range("B1").Value = "No Induk"
$ws.Range("C1").Value = "NISN"
$ws.Range("D1").Value = "JK"
$ws.Range("E1").Value = "Tempat Lahir"
$ws.Range("F1").Value = "Tanggal Lahir"
$ws.Range("G1").Value = "Agama"
$ws.Range("H1").Value = "Status dlm Keluarga"
$ws.Range("I1").Value = "Anak ke-"
$ws.Range("J1").Value = "Alamat"
$ws.Range("K1").Value = "Asal Sekolah"
$ws.Range("L1").Value = "Diterima dikelas"
$ws.Range("M1").Value = "Tgl diterima"
$ws.Range("N1").Value = "Ayah"
$ws.Range("O1").Value = "Pekerjaan Ayah"
$ws.Range("P1").Value = "Ibu"
$ws.Range("Q1").Value = "Pekerjaan Ibu"
$ws.Range("R1").Value = "Wali"
$ws.Range("S1").Value = "Pekerjaan Wali"

# --- 3. Row 2 ----------------------------------------------------------------
$ws.Range("A2").Value = "maman abdurahman"
$ws.Range("B2").Value = 12345
$ws.Range("C2").Value = 2342
$ws.Range("D2").Value = "M"
$ws.Range("E2").Value = "bogr"
Set-TextValue $ws.Range("F2") "2020-08-04"
$ws.Range("G2").Value = "adf"
$ws.Range("H2").Value = "dfadf"
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = "afaf"
$ws.Range("K2").Value = "adf"
$ws.Range("L2").Value = "dfsa"
Set-TextValue $ws.Range("M2") "2020-08-14"
$ws.Range("N2").Value = "dfafd"
$ws.Range("O2").Value = "dfdfsf"
$ws.Range("P2").Value = "adf"
$ws.Range("Q2").Value = "adf"
$ws.Range("R2").Value = "adf"
$ws.Range("S2").Value = "adfas"

# --- 4. Row 3 ----------------------------------------------------------------
$ws.Range("A3").Value = "Aisayah "
$ws.Range("B3").Value = 2345
$ws.Range("C3").Value = 134
$ws.Range("D3").Value = "M"
$ws.Range("E3").Value = "adsf"
Set-TextValue $ws.Range("F3") "2020-08-02"
$ws.Range("G3").Value = "df"
$ws.Range("H3").Value = "adfs"
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = "dfasf"
$ws.Range("K3").Value = "adfa"
$ws.Range("L3").Value = "aa"
Set-TextValue $ws.Range("M3") "2020-08-02"
$ws.Range("N3").Value = "adfaf"
$ws.Range("O3").Value = "afasf"
$ws.Range("P3").Value = "afdadf"
$ws.Range("Q3").Value = "dafas"
$ws.Range("R3").Value = "adfas"
$ws.Range("S3").Value = "sadfa"

# --- 5. Selection state ------------------------------------------------------
$ws.Range("A1:A1048576").Select()
